# Singapore Premier League workbook update (12-06-2024 23:38)
#
# This edit does two things:
#  1. Swaps the team names "Lion City Sailors FC" and "DPMM FC" wherever they
#     are used as a HomeTeam/AwayTeam value on the results sheet (the two
#     teams were mislabelled against each other in the source feed).
#  2. A handful of fixture rows had been matched to the wrong match id/odds
#     set; for those row pairs the entire data row (id, teams, score, odds,
#     …) is swapped between the two rows, leaving only the Div and Date
#     columns (which are correct and identical for the pair) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$used = $ws.UsedRange
$lastRow = $used.Row + $used.Rows.Count - 1

# --- Step 1: swap the "Lion City Sailors FC" / "DPMM FC" team names -------
# Teams only ever show up in columns E (HomeTeam) and F (AwayTeam); use a
# placeholder so the two values trade places rather than collapsing to one.
$placeholder = "__TEAM_SWAP_PLACEHOLDER__"

foreach ($colLetter in @("E", "F")) {
    $col = $ws.Range($colLetter + "1").Column
    for ($r = 2; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, $col)
        if ($cell.Value2 -eq "Lion City Sailors FC") {
            $cell.Value2 = $placeholder
        }
    }
    for ($r = 2; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, $col)
        if ($cell.Value2 -eq "DPMM FC") {
            $cell.Value2 = "Lion City Sailors FC"
        }
    }
    for ($r = 2; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, $col)
        if ($cell.Value2 -eq $placeholder) {
            $cell.Value2 = "DPMM FC"
        }
    }
}

# --- Step 2: swap mismatched fixture rows ----------------------------------
# Columns: B=id C=Div D=Date E=HomeTeam F=AwayTeam G..AD = score/odds data.
# Div (C) and Date (D) are shared/correct for each pair and stay put; every
# other column swaps between the two rows.
$rowPairs = @(
    @(6, 7),
    @(18, 19),
    @(20, 21),
    @(43, 44),
    @(54, 55)
)

# Explicit column list B, E, F, G, H, ..., Z, AA, AB, AC, AD (skip A, C, D)
$colLetters = @("B", "E", "F", "G", "H", "I", "J", "K", "L", "M", "N", "O", "P", "Q", "R", "S", "T", "U", "V", "W", "X", "Y", "Z", "AA", "AB", "AC", "AD")

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    foreach ($colLetter in $colLetters) {
        $col = $ws.Range($colLetter + "1").Column
        $cell1 = $ws.Cells.Item($r1, $col)
        $cell2 = $ws.Cells.Item($r2, $col)
        $v1 = $cell1.Value2
        $v2 = $cell2.Value2
        $cell1.Value2 = $v2
        $cell2.Value2 = $v1
    }
}
